{"js": "// Replace each three-digit-divided-by-one-digit expression with its updated value.\n// Each expression text is unique in the document, so a simple search+replace per pair is unambiguous.\nconst replacements = [\n  [\"148\u00f73=\", \"200\u00f72=\"],\n  [\"640\u00f73=\", \"567\u00f78=\"],\n  [\"281\u00f73=\", \"496\u00f75=\"],\n  [\"670\u00f73=\", \"948\u00f73=\"],\n  [\"435\u00f72=\", \"723\u00f75=\"],\n  [\"862\u00f72=\", \"749\u00f73=\"],\n  [\"167\u00f72=\", \"560\u00f72=\"],\n  [\"685\u00f77=\", \"957\u00f73=\"],\n  [\"944\u00f79=\", \"982\u00f74=\"],\n  [\"692\u00f78=\", \"266\u00f76=\"],\n  [\"606\u00f76=\", \"662\u00f74=\"],\n  [\"836\u00f72=\", \"375\u00f75=\"],\n  [\"442\u00f79=\", \"545\u00f72=\"],\n  [\"833\u00f74=\", \"984\u00f78=\"],\n  [\"257\u00f72=\", \"350\u00f78=\"],\n  [\"489\u00f75=\", \"739\u00f74=\"],\n  [\"737\u00f78=\", \"374\u00f78=\"],\n  [\"771\u00f72=\", \"930\u00f72=\"],\n  [\"433\u00f76=\", \"687\u00f73=\"],\n  [\"994\u00f76=\", \"847\u00f72=\"],\n  [\"838\u00f73=\", \"651\u00f73=\"],\n  [\"485\u00f73=\", \"218\u00f79=\"],\n  [\"225\u00f75=\", \"949\u00f74=\"],\n  [\"908\u00f74=\", \"924\u00f73=\"],\n  [\"165\u00f79=\", \"985\u00f74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit \u00f7 one-digit expression to its new value.\n# Every expression text is unique within the document, so Find/Replace\n# on the exact text is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ old = \"148\u00f73=\"; new = \"200\u00f72=\" }\n    @{ old = \"640\u00f73=\"; new = \"567\u00f78=\" }\n    @{ old = \"281\u00f73=\"; new = \"496\u00f75=\" }\n    @{ old = \"670\u00f73=\"; new = \"948\u00f73=\" }\n    @{ old = \"435\u00f72=\"; new = \"723\u00f75=\" }\n    @{ old = \"862\u00f72=\"; new = \"749\u00f73=\" }\n    @{ old = \"167\u00f72=\"; new = \"560\u00f72=\" }\n    @{ old = \"685\u00f77=\"; new = \"957\u00f73=\" }\n    @{ old = \"944\u00f79=\"; new = \"982\u00f74=\" }\n    @{ old = \"692\u00f78=\"; new = \"266\u00f76=\" }\n    @{ old = \"606\u00f76=\"; new = \"662\u00f74=\" }\n    @{ old = \"836\u00f72=\"; new = \"375\u00f75=\" }\n    @{ old = \"442\u00f79=\"; new = \"545\u00f72=\" }\n    @{ old = \"833\u00f74=\"; new = \"984\u00f78=\" }\n    @{ old = \"257\u00f72=\"; new = \"350\u00f78=\" }\n    @{ old = \"489\u00f75=\"; new = \"739\u00f74=\" }\n    @{ old = \"737\u00f78=\"; new = \"374\u00f78=\" }\n    @{ old = \"771\u00f72=\"; new = \"930\u00f72=\" }\n    @{ old = \"433\u00f76=\"; new = \"687\u00f73=\" }\n    @{ old = \"994\u00f76=\"; new = \"847\u00f72=\" }\n    @{ old = \"838\u00f73=\"; new = \"651\u00f73=\" }\n    @{ old = \"485\u00f73=\"; new = \"218\u00f79=\" }\n    @{ old = \"225\u00f75=\"; new = \"949\u00f74=\" }\n    @{ old = \"908\u00f74=\"; new = \"924\u00f73=\" }\n    @{ old = \"165\u00f79=\"; new = \"985\u00f74=\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
